# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to "_FV2410" / "_FV2504"
# - Freeze the header row
# - Wrap the data range in an Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Header row renames ------------------------------------------------
# Columns A:J carried the "_old" suffix, column K is "diff", columns L:U
# carried the "_new" suffix. Swap the suffixes for the new release labels.
$oldHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$newHeaders = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap A1:U57 in an Excel Table --------------------------------------
$tableRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$tbl.Name = "Table1"
